# feat: add 2022-Q3 data
#
# The workbook has a "总计" (totals) sheet plus a single quarter sheet
# named "2022-Q1". We duplicate the "2022-Q1" sheet (preserving its data
# unchanged) to a new sheet placed right after it, rename the ORIGINAL
# sheet to "2022-Q3" and overwrite its contents with the new quarter's
# fund-holding data, then add the new quarter's summary row to "总计".

$wb = $excel.ActiveWorkbook

# ---- 1. Locate existing sheets -------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)   # "总计"
$q1Sheet    = $wb.Worksheets.Item(2)   # currently named "2022-Q1"

# ---- 2. Duplicate the Q1 sheet so its data survives under its own tab ----------
$q1Sheet.Copy($null, $q1Sheet)
$q1Copy = $wb.Worksheets.Item(3)

# Rename the original first (frees up the "2022-Q1" name), then rename the
# freshly-created copy to take over that name.
$q1Sheet.Name = "2022-Q3"
$q1Copy.Name  = "2022-Q1"

$q3Sheet = $q1Sheet

# ---- 3. Re-style the Q3 sheet's header row + A2 to match "总计"'s accent style --
$totalSheet.Range("B1").Copy()
$q3Sheet.Range("B1:H1").PasteSpecial(-4122)
$totalSheet.Range("A2").Copy()
$q3Sheet.Range("A2").PasteSpecial(-4122)

# ---- 4. Overwrite the Q3 sheet's fund row with the new quarter's data ----------
# B2 ("012995") and D2:G2 (decimal-looking numbers) must stay literal TEXT
# (matching the source data's inlineStr cells), so force Text format before
# assigning them and drop back to the default "Normal" style afterwards
# (Style="Normal" clears the number-format override without leaving any
# `s="..."` attribute behind, since it resolves back to the workbook's
# built-in style 0).
$textCells = @("B2", "D2", "E2", "F2", "G2")
foreach ($addr in $textCells) {
    $q3Sheet.Range($addr).NumberFormat = "@"
}
$q3Sheet.Range("B2").Value = "012995"
$q3Sheet.Range("C2").Value = "嘉实策略视野三年持有期混合"
$q3Sheet.Range("D2").Value = "3.04"
$q3Sheet.Range("E2").Value = "64.76"
$q3Sheet.Range("F2").Value = "2.47"
$q3Sheet.Range("G2").Value = "0.0751"
foreach ($addr in $textCells) {
    $q3Sheet.Range($addr).Style = "Normal"
}

$q3Sheet.Range("H2").Value = 10

# ---- 5. Update the "总计" summary sheet -----------------------------------------
# Row 2 becomes the new Q3 totals; the old Q1 totals move down to row 3.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q1"
$totalSheet.Range("C3").Value = 1
$totalSheet.Range("D3").Value = 0.48

$totalSheet.Range("A2").Copy()
$totalSheet.Range("A3").PasteSpecial(-4122)

$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("D2").Value = 0.08
